$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '68.341.53'
$ws.Range("E2").Value = '  +2.90%  '

# Row 3
$ws.Range("D3").Value = '2.539.50'
$ws.Range("E3").Value = '  +1.43%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.83'
$ws.Range("E5").Value = '  +2.27%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.87'
$ws.Range("E6").Value = '  +3.82%  '

# Row 7
$ws.Range("E7").Value = '  -0.06%  '

# Row 8
$ws.Range("E8").Value = '  +1.56%  '

# Row 9
$ws.Range("D9").Value = '2.538.46'
$ws.Range("E9").Value = '  +1.44%  '

# Row 10
$ws.Range("E10").Value = '  +2.69%  '

# Row 11
$ws.Range("E11").Value = '  +2.96%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.15'
$ws.Range("E12").Value = '  +0.80%  '

# Row 13
$ws.Range("E13").Value = '  -0.26%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.03'
$ws.Range("E14").Value = '  +1.61%  '

# Row 15
$ws.Range("D15").Value = '3.002.20'
$ws.Range("E15").Value = '  +1.52%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  +2.67%  '

# Row 17
$ws.Range("D17").Value = '68.279.17'
$ws.Range("E17").Value = '  +2.96%  '

# Row 18
$ws.Range("D18").Value = '2.543.05'
$ws.Range("E18").Value = '  +1.02%  '

# Row 19
$ws.Range("E19").Value = '  +4.83%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.56'
$ws.Range("E20").Value = '  +3.09%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '369.41'
$ws.Range("E21").Value = '  +6.49%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.22'
$ws.Range("E22").Value = '  +0.90%  '

# Row 23
$ws.Range("E23").Value = '  +2.52%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.95'
$ws.Range("E24").Value = '  -0.75%  '

# Row 25
$ws.Range("E25").Value = '  -0.03%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.33'
$ws.Range("E26").Value = '  +3.49%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.01'
$ws.Range("E27").Value = '  +2.15%  '

# Row 28
$ws.Range("D28").Value = '2.672.34'
$ws.Range("E28").Value = '  +1.73%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.995'
$ws.Range("E29").Value = '  -0.52%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0999'
$ws.Range("E30").Value = '  +2.62%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '541.94'
$ws.Range("E31").Value = '  +2.85%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.31'
$ws.Range("E32").Value = '  +3.07%  '

# Row 33
$ws.Range("E33").Value = '  +2.39%  '

# Row 34
$ws.Range("E34").Value = '  +2.97%  '

# Row 35
$ws.Range("E35").Value = '  -0.23%  '

# Row 36
$ws.Range("E36").Value = '  +0.01%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.47'
$ws.Range("E37").Value = '  +1.46%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.50'
$ws.Range("E38").Value = '  +0.80%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.94'
$ws.Range("E39").Value = '  +1.87%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.73'
$ws.Range("E40").Value = '  +2.02%  '

# Row 41
$ws.Range("E41").Value = '  +2.22%  '

# Row 42
$ws.Range("B42").Value = 'PolygonEcosystemToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.357'
$ws.Range("E42").Value = '  +0.91%  '

# Row 43
$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.24'
$ws.Range("E43").Value = '  +3.16%  '

# Row 44
$ws.Range("E44").Value = '  +2.70%  '

# Row 45
$ws.Range("E45").Value = '  +0.06%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.565'
$ws.Range("E46").Value = '  +1.69%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '147.79'
$ws.Range("E47").Value = '  -0.26%  '

# Row 48
$ws.Range("D48").Value = '0.0₆0281'
$ws.Range("E48").Value = '  +3.40%  '

# Row 49
$ws.Range("E49").Value = '  +2.13%  '

# Row 50
$ws.Range("E50").Value = '  +0.20%  '

# Row 51
$ws.Range("E51").Value = '  +0.97%  '
